$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename "Perfil V-Slot 2020" to "Perfil V-Slot 2040"
$ws.Range("A5").Value = "Perfil V-Slot 2040 "

# 2. Update row 16 (Turecas tipo martillo M4): quantity 6 -> 8, add a note
$ws.Range("B16").Value = 8
$ws.Range("F16").Value = "Para sujeción de correas y finales de carrera"

# 3. Insert two new rows after row 24 (before the old row 26) for the new
#    "Tornillos M4x10" / "Tornillos M2x8" items, pushing everything below
#    down by two rows.
$ws.Rows("25:26").Insert()

# Copy the formatting (style) used on row 17 onto the new rows/cells so the
# new cells match the rest of the table's look (font/style index), without
# touching the untouched D/E columns.
$ws.Range("A17").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("C17").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("F17").Copy()
$ws.Range("F25").PasteSpecial(-4122)

$ws.Range("A17").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("C17").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("F17").Copy()
$ws.Range("F26").PasteSpecial(-4122)

$ws.Range("A25").Value = "Tornillos M4x10"
$ws.Range("B25").Value = 8
$ws.Range("C25").Value = "10mm"
$ws.Range("F25").Value = "Para sujeción de correas y finales de carrera"

$ws.Range("A26").Value = "Tornillos M2x8"
$ws.Range("B26").Value = 4
$ws.Range("C26").Value = "8mm"
$ws.Range("F26").Value = "Para los finales de carrera"

$ws.Range("A1").Select()

# 4. Fill in the previously-blank row 33 with a new "finales de carrera" item.
$ws.Range("A33").Value = "finales de carrera"
$ws.Range("B33").Value = 2
$ws.Range("D33").Value = "https://amzn.to/2vvWAP3"

$ws.Range("E33").Select()
